$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve textual formatting for numeric-looking price values (avoid Excel auto-number conversion)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

# Apply updated Price (D) and Volume(1h) (E) values
$ws.Range("D2").Value = '67.594.42'
$ws.Range("E2").Value = '  +4.29%  '
$ws.Range("D3").Value = '3.253.76'
$ws.Range("E3").Value = '  +3.42%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '578.30'
$ws.Range("E5").Value = '  +2.28%  '
$ws.Range("D6").Value = '181.75'
$ws.Range("E6").Value = '  +6.25%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '0.592'
$ws.Range("E8").Value = '  -4.26%  '
$ws.Range("D9").Value = '3.255.36'
$ws.Range("E9").Value = '  +3.67%  '
$ws.Range("E10").Value = '  +4.45%  '
$ws.Range("E11").Value = '  +3.41%  '
$ws.Range("D12").Value = '0.414'
$ws.Range("E12").Value = '  +4.89%  '
$ws.Range("D13").Value = '3.825.09'
$ws.Range("E13").Value = '  +3.70%  '
$ws.Range("E14").Value = '  +1.65%  '
$ws.Range("D15").Value = '28.17'
$ws.Range("E15").Value = '  +3.71%  '
$ws.Range("D16").Value = '67.567.23'
$ws.Range("E16").Value = '  +4.41%  '
$ws.Range("E17").Value = '  +3.00%  '
$ws.Range("D18").Value = '3.258.87'
$ws.Range("E18").Value = '  +3.61%  '
$ws.Range("E19").Value = '  +1.91%  '
$ws.Range("D20").Value = '13.51'
$ws.Range("E20").Value = '  +5.02%  '
$ws.Range("D21").Value = '376.49'
$ws.Range("E21").Value = '  +5.90%  '
$ws.Range("E22").Value = '  +5.33%  '
$ws.Range("E23").Value = '  -0.24%  '
$ws.Range("D24").Value = '71.28'
$ws.Range("E24").Value = '  +4.08%  '
$ws.Range("D25").Value = '0.511'
$ws.Range("E25").Value = '  +2.33%  '
$ws.Range("D26").Value = '0.0000119'
$ws.Range("E26").Value = '  +1.14%  '
$ws.Range("D27").Value = '9.64'
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("E28").Value = '  +3.48%  '
$ws.Range("D29").Value = '1.01'
$ws.Range("E29").Value = '  +0.51%  '
$ws.Range("E30").Value = '  +6.47%  '
$ws.Range("E31").Value = '  +3.96%  '
$ws.Range("D32").Value = '22.65'
$ws.Range("E32").Value = '  +3.04%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("E34").Value = '  +6.03%  '
$ws.Range("D35").Value = '6.88'
$ws.Range("E35").Value = '  +3.61%  '
$ws.Range("D36").Value = '162.43'
$ws.Range("E36").Value = '  +6.08%  '
$ws.Range("E37").Value = '  +3.54%  '
$ws.Range("D38").Value = '0.854'
$ws.Range("E38").Value = '  +3.07%  '
$ws.Range("E39").Value = '  +6.68%  '
$ws.Range("D40").Value = '6.82'
$ws.Range("E40").Value = '  +13.16%  '
$ws.Range("D41").Value = '26.83'
$ws.Range("E41").Value = '  +3.18%  '
$ws.Range("E42").Value = '  +2.96%  '
$ws.Range("D43").Value = '362.82'
$ws.Range("E43").Value = '  +13.58%  '
$ws.Range("D44").Value = '4.47'
$ws.Range("E44").Value = '  +7.05%  '
$ws.Range("D45").Value = '2.746.75'
$ws.Range("E45").Value = '  +3.41%  '
$ws.Range("D46").Value = '25.52'
$ws.Range("E46").Value = '  +5.49%  '
$ws.Range("D47").Value = '40.68'
$ws.Range("E47").Value = '  +4.28%  '
$ws.Range("D48").Value = '0.0674'
$ws.Range("E48").Value = '  +3.05%  '
$ws.Range("E49").Value = '  +2.25%  '
$ws.Range("E50").Value = '  +6.84%  '
$ws.Range("E51").Value = '  +0.16%  '
